# Working language spreadsheet integration.
# Updates font filenames (tff -> ttf typo fix, and Korean text font switched
# to a shared "defaultFont"), renames the app name, and updates the
# active selection on the dialogue sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dialogue")

# Fix font file name typos (.tff -> .ttf) and repoint the Korean text font
# to a shared default font entry.
$ws.Range("B2").Value = "heading_EN.ttf"
$ws.Range("C2").Value = "heading_EN.ttf"
$ws.Range("D2").Value = "heading_KOR.ttf"

$ws.Range("B3").Value = "text_EN.ttf"
$ws.Range("C3").Value = "text_CZ.ttf"
$ws.Range("D3").Value = "defaultFont"

# Update the app display name.
$ws.Range("B5").Value = "Shooting Stars :3"

# Move the active selection / view to D3, matching the frozen pane scrolled
# so column C is the left-most visible column after the split.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D3").Select()
